# Updates to the "想去人数" (F column) counts on the 展览 and 全部类型 sheets,
# reflecting refreshed data from a gh-pages generation run.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 829
$ws1.Range("F6").Value = 3750
$ws1.Range("F7").Value = 2500
$ws1.Range("F9").Value = 2329
$ws1.Range("F15").Value = 87
$ws1.Range("F21").Value = 421
$ws1.Range("F29").Value = 30
$ws1.Range("F31").Value = 776
$ws1.Range("F32").Value = 801
$ws1.Range("F33").Value = 1895
$ws1.Range("F37").Value = 548
$ws1.Range("F38").Value = 1182
$ws1.Range("F40").Value = 397

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 829
$ws4.Range("F6").Value = 3750
$ws4.Range("F7").Value = 2500
$ws4.Range("F9").Value = 2329
$ws4.Range("F15").Value = 87
$ws4.Range("F21").Value = 421
$ws4.Range("F32").Value = 30
$ws4.Range("F34").Value = 776
$ws4.Range("F36").Value = 801
$ws4.Range("F37").Value = 1895
$ws4.Range("F44").Value = 548
$ws4.Range("F45").Value = 1182
$ws4.Range("F47").Value = 397
